$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the (now-removed) speaker names in column D with the new
# placeholder values that were typed into the shared-string table.
# "dds" is entered before "dsd" so the new shared-string entries land
# in the same order as in the saved workbook.
$ws.Range("D6").Value = "dds"
$ws.Range("D4").Value = "dsd"
$ws.Range("D5").Value = "dsd"
$ws.Range("D8").Value = "dsd"
$ws.Range("D10").Value = "dsd"

# Update the start/end times for the last two agenda rows (moved from
# the afternoon into the evening).
$ws.Range("A17").Value = 0.91666666666666663
$ws.Range("B17").Value = 0.95833333333333337
$ws.Range("A18").Value = 0.95833333333333337
$ws.Range("B18").Value = 0.97916666666666663

# Re-apply the original time-cell formatting (the plain Value write
# above resets the quote-prefixed number-format style), by copying the
# format from the still-untouched row above.
$ws.Range("A16:B16").Copy()
$ws.Range("A17:B18").PasteSpecial(-4122)

# Update the active selection to match the saved worksheet view.
$ws.Range("B19").Select()
